# Updated cryptos list — applies Price (D) and Volume(1h) (E) cell text updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.313.38'
$ws.Range('E2').Value = '  -2.49%  '
$ws.Range('D3').Value = '1.708.71'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '223.56'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5308'
$ws.Range('E6').Value = '  -2.69%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2655'
$ws.Range('E8').Value = '  -4.84%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06599'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.92'
$ws.Range('E10').Value = '  -4.47%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07650'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.583'
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('D13').Value = '1.719.21'
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('D14').Value = '1.944.70'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5743'
$ws.Range('E15').Value = '  -5.04%  '
$ws.Range('D16').Value = '0.0₅8193'
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.62'
$ws.Range('E17').Value = '  -3.11%  '
$ws.Range('D18').Value = '27.309.08'
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '216.54'
$ws.Range('E19').Value = '  -4.46%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.42'
$ws.Range('E22').Value = '  -5.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.973'
$ws.Range('E23').Value = '  -4.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '142.39'
$ws.Range('E25').Value = '  -3.15%  '
$ws.Range('E26').Value = '  +6.78%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1216'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.266'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '16.35'
$ws.Range('E29').Value = '  -4.97%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05386'
$ws.Range('E30').Value = '  -4.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.293'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.508'
$ws.Range('E32').Value = '  -5.56%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.425'
$ws.Range('E33').Value = '  -3.56%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.643'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.423'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9479'
$ws.Range('E37').Value = '  -3.95%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5861'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01633'
$ws.Range('E39').Value = '  -3.14%  '
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.004'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').Value = '1.042.74'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8422'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '101.14'
$ws.Range('E44').Value = '  -1.28%  '
$ws.Range('D45').Value = '1.850.82'
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  -1.51%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '58.12'
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('E48').Value = '  +1.74%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.119'
$ws.Range('E50').Value = '  -2.75%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06506'
$ws.Range('E51').Value = '  +9.38%  '
